$wb = $excel.ActiveWorkbook

# --- Rename the 8 table sheets (new random suffixes), keep sheetId/order the same ---
$wb.Worksheets.Item("Unnamed_Table_2ce14a3c_6").Name  = "Unnamed_Table_f3ac1cfe_6"
$wb.Worksheets.Item("Unnamed_Table_023cf8a7_11").Name = "Unnamed_Table_a38c5d5c_11"
$wb.Worksheets.Item("Unnamed_Table_52166d9a_11").Name = "Unnamed_Table_f6dd28ed_11"
$wb.Worksheets.Item("Unnamed_Table_c4f900c3_11").Name = "Unnamed_Table_617355eb_11"
$wb.Worksheets.Item("Unnamed_Table_2c9cf0c1_11").Name = "Unnamed_Table_42baba9f_11"
$wb.Worksheets.Item("Unnamed_Table_c37fd4ca_12").Name = "Unnamed_Table_1122d327_12"
$wb.Worksheets.Item("Unnamed_Table_9107e2cf_13").Name = "Unnamed_Table_9ccab279_13"
$wb.Worksheets.Item("Unnamed_Table_a279fd4c_13").Name = "Unnamed_Table_20e3b24b_13"

# --- Update the Summary sheet's table_name / full_name columns to match ---
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("A2").Value = "Unnamed_Table_f3ac1cfe_6"

$summary.Range("A3").Value = "Unnamed_Table_a38c5d5c_11"
$summary.Range("B3").Value = "OSS vs Proprietary Scores by Language Model"

$summary.Range("A4").Value = "Unnamed_Table_f6dd28ed_11"
$summary.Range("B4").Value = "Language Model Performance by License Category"

$summary.Range("A5").Value = "Unnamed_Table_617355eb_11"
$summary.Range("B5").Value = "C# Code Generation Metric Comparison"

$summary.Range("A6").Value = "Unnamed_Table_42baba9f_11"
$summary.Range("B6").Value = "C Sharp Code Generation Metrics Comparison"

$summary.Range("A7").Value = "Unnamed_Table_1122d327_12"

$summary.Range("A8").Value = "Unnamed_Table_9ccab279_13"
$summary.Range("B8").Value = "C Sharp Code Davinci 002 BLEU Scores by Shot Source"

$summary.Range("A9").Value = "Unnamed_Table_20e3b24b_13"
$summary.Range("B9").Value = "LLM Accuracy on Identifiers and Non identifiers"

# --- Fix header text on the individual table sheets ---
$wb.Worksheets.Item("Unnamed_Table_a38c5d5c_11").Range("E1").Value = "Category-p-value"
$wb.Worksheets.Item("Unnamed_Table_20e3b24b_13").Range("A1").Value = "Language"
